# RPA datasets push 2023-12-14
# The IPO price-discovery table had two rows (하나스팩30호 / 디에스단석) whose
# "확정공모가" (confirmed offering price) was still pending ("-") when the
# sheet was first captured. Both IPOs have since priced, so fill in the
# now-known confirmed offering prices.
#
# These text values look numeric, so they must be forced to remain text
# (matching how the rest of this column is stored) instead of being
# auto-coerced into numeric cells by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2000"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "100000"
$ws.Range("D7").Style = "Normal"
